# Generate Report for Archive
#
# 1. Change status text "Ready for handoff" -> "In Translation" on all three sheets
#    (Overview columns E/F, zh-cn column C, de-de column C).
# 2. Narrow the "Status" column width(s) (Overview: E & F, zh-cn: C, de-de: C).
#    The stored column width is quantized by the engine to steps of 1/6 of a
#    character; requesting 12.5 lands on the nearest achievable width to the
#    target 13.4101845877511 (i.e. 13.333333333333334 = 80/6).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColumnWidth = 12.5

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
